$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 104.00990204435887
$ws.Range("C2").Value = 104.23615047329335
$ws.Range("D2").Value = 106.13121722688041
$ws.Range("E2").Value = 104.92882886823251

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 104.46470509868946
$ws.Range("C3").Value = 105.5223897173454
$ws.Range("D3").Value = 105.16444695554922
$ws.Range("E3").Value = 104.69558980577989

# Update selection to reflect the new active range
$ws.Range("B1:E3").Select() | Out-Null
